# Updated MCH102 to MCH251: populate the (previously header-only) sheet
# with its first data row, describing collection MCH131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values --------------------------------------------------
# (column B / alternativeIdentifiers, D / date_s and H / file_path
# stay blank for this record)
$ws.Range("A2").Value = "MCH131"
$ws.Range("C2").Value = "JEROME MAALE PAPERS, EXZAM RESULTS, THE AFRICAN ROCK OF FREEDOM"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION:  CABINET 1B | GRAP COUNT NUMER: NONE"

# D2 (date_s) and H2 (file_path) are intentionally left blank, but still
# carry the new row's cell formatting like the rest of the row.
$rngA = $ws.Range("A2")
$rngA.Font.Name = "Calibri"
$rngA.Font.Size = 10
$rngA.Font.ThemeColor = 1

$rngRest = $ws.Range("C2:H2")
$rngRest.Font.Name = "Calibri"
$rngRest.Font.Size = 10
$rngRest.Font.ThemeColor = 1

# --- View state ------------------------------------------------------
# Re-establish the frozen header row/pane and move the active selection
# down onto the newly added data row.
$ws.Range("A2").Select() | Out-Null
$win = $excel.ActiveWindow
$win.FreezePanes = $true
$ws.Range("A4").Select() | Out-Null
